# Trading update: 2026-02-18 10:19:54
# - Trade #7 (row 8 on "All Trades") closes out its "latest trade" snapshot
#   fields (Capital After / slippage / confidence / entry reason / duration)
#   since it is no longer the most recent trade, and gets an Exit Price of 0.
# - Three new trades (#8, #9, #10) are appended as rows 9-11.
# - Trade #10 is the new "latest trade", so it carries the snapshot fields
#   that used to live on trade #7's row.
# - The per-strategy "MarketMaking" sheet's row 2 (latest trade for that
#   strategy) is refreshed to mirror trade #10.

$wb = $excel.ActiveWorkbook
$all = $wb.Worksheets.Item("All Trades")
$mm = $wb.Worksheets.Item("MarketMaking")

# Date/Time columns hold plain text like "2026-02-18" / "10:18:32" - force
# text formatting first so COM does not reinterpret them as date/time
# serials.
$all.Range("B9:C11").NumberFormat = "@"
$mm.Range("B2:C2").NumberFormat = "@"

# --- "All Trades": update row 8 (trade #7) -------------------------------
# It is no longer the latest trade, so exit price becomes a real 0 and the
# "latest snapshot" columns (K:Q) are cleared.
$all.Cells.Item(8, 7).Value = 0          # G8 Exit Price
$all.Cells.Item(8, 11).Value = ""        # K8 Capital After
$all.Cells.Item(8, 12).Value = ""        # L8 Entry Slippage (bps)
$all.Cells.Item(8, 13).Value = ""        # M8 Exit Slippage (bps)
$all.Cells.Item(8, 14).Value = ""        # N8 Confidence
$all.Cells.Item(8, 15).Value = ""        # O8 Entry Reason
$all.Cells.Item(8, 16).Value = ""        # P8 Exit Reason (already blank)
$all.Cells.Item(8, 17).Value = ""        # Q8 Duration (min)

# --- "All Trades": append trade #8 as row 9 ------------------------------
$all.Cells.Item(9, 1).Value = 8
$all.Cells.Item(9, 2).Value = "2026-02-18"
$all.Cells.Item(9, 3).Value = "10:18:32"
$all.Cells.Item(9, 4).Value = "MarketMaking"
$all.Cells.Item(9, 5).Value = "DOWN"
$all.Cells.Item(9, 6).Value = 0.18
$all.Cells.Item(9, 7).Value = 0
$all.Cells.Item(9, 8).Value = "OPEN"
$all.Cells.Item(9, 9).Value = 0
$all.Cells.Item(9, 10).Value = 0

# --- "All Trades": append trade #9 as row 10 -----------------------------
$all.Cells.Item(10, 1).Value = 9
$all.Cells.Item(10, 2).Value = "2026-02-18"
$all.Cells.Item(10, 3).Value = "10:18:38"
$all.Cells.Item(10, 4).Value = "MarketMaking"
$all.Cells.Item(10, 5).Value = "UP"
$all.Cells.Item(10, 6).Value = 0.9
$all.Cells.Item(10, 7).Value = 0
$all.Cells.Item(10, 8).Value = "OPEN"
$all.Cells.Item(10, 9).Value = 0
$all.Cells.Item(10, 10).Value = 0

# --- "All Trades": append trade #10 as row 11 (new latest trade) --------
$all.Cells.Item(11, 1).Value = 10
$all.Cells.Item(11, 2).Value = "2026-02-18"
$all.Cells.Item(11, 3).Value = "10:19:39"
$all.Cells.Item(11, 4).Value = "MarketMaking"
$all.Cells.Item(11, 5).Value = "DOWN"
$all.Cells.Item(11, 6).Value = 0.57
$all.Cells.Item(11, 8).Value = "OPEN"
$all.Cells.Item(11, 9).Value = 0
$all.Cells.Item(11, 10).Value = 0
$all.Cells.Item(11, 11).Value = 100
$all.Cells.Item(11, 12).Value = 0
$all.Cells.Item(11, 13).Value = 0
$all.Cells.Item(11, 14).Value = 0.6
$all.Cells.Item(11, 15).Value = "Normal spread capture: 619 bps"
$all.Cells.Item(11, 17).Value = 0

# --- "MarketMaking": refresh row 2 (latest MarketMaking trade = #10) ----
$mm.Cells.Item(2, 1).Value = 10
$mm.Cells.Item(2, 3).Value = "10:19:39"
$mm.Cells.Item(2, 5).Value = "DOWN"
$mm.Cells.Item(2, 6).Value = 0.57
$mm.Cells.Item(2, 15).Value = "Normal spread capture: 619 bps"

# The values are already cached as plain text at this point, so the
# temporary "@" number format can be reset back to the default "Normal"
# style - this keeps the written cells style-free, matching every other
# cell in these sheets (which all use the workbook's default style).
$all.Range("B9:C11").Style = "Normal"
$mm.Range("B2:C2").Style = "Normal"
